$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data grid is A1:L12 (A/1 are labels). Numeric payload is B2:L12.
# Every numeric cell whose current value is not exactly 1 gets transformed
# via new = old * 2 + 1. Cells with value 1 (diagonal / symmetric identity
# entries) and non-numeric (inline string / blank) cells are left untouched.
for ($row = 2; $row -le 12; $row++) {
    for ($col = 2; $col -le 12; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value2
        if ($val -is [double] -and $val -ne 1) {
            $cell.Value2 = $val * 2 + 1
        }
    }
}
